$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.49280000000001
$ws.Range("A10").Value = -20.46099999999997
$ws.Range("A12").Value = -22.40910000000003
$ws.Range("D13").Value = -7.651900000000001
$ws.Range("A18").Value = -22.32020000000003
$ws.Range("E20").Value = 12.2112
